$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.817.64"
$ws.Range("E2").Value = "  +3.04%  "
$ws.Range("D3").Value = "3.975.46"
$ws.Range("E3").Value = "  +0.96%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.09"
$ws.Range("E5").Value = "  +9.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.03"
$ws.Range("E6").Value = "  +8.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.685"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.748"
$ws.Range("E9").Value = "  +1.84%  "
$ws.Range("E10").Value = "  +1.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.10"
$ws.Range("E11").Value = "  -0.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000318"
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.86"
$ws.Range("E13").Value = "  +3.04%  "
$ws.Range("D14").Value = "4.607.37"
$ws.Range("E14").Value = "  +1.07%  "
$ws.Range("D15").Value = "3.983.53"
$ws.Range("E15").Value = "  +1.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.27"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.04"
$ws.Range("E17").Value = "  +2.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.33"
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("D20").Value = "72.633.93"
$ws.Range("E20").Value = "  +2.97%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "434.48"
$ws.Range("E21").Value = "  +2.80%  "
$ws.Range("E22").Value = "  +13.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "96.01"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("E24").Value = "  -3.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.28"
$ws.Range("E25").Value = "  +1.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.46"
$ws.Range("E26").Value = "  +22.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.17"
$ws.Range("E27").Value = "  -1.53%  "
$ws.Range("E28").Value = "  +0.34%  "
$ws.Range("E29").Value = "  +1.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.30"
$ws.Range("E30").Value = "  +0.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.83"
$ws.Range("E31").Value = "  +2.84%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.66"
$ws.Range("E32").Value = "  +3.34%  "
$ws.Range("E33").Value = "  +2.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "681.88"
$ws.Range("E34").Value = "  -0.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "48.34"
$ws.Range("E35").Value = "  -2.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "69.62"
$ws.Range("E36").Value = "  +8.74%  "
$ws.Range("E37").Value = "  +7.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.434"
$ws.Range("E38").Value = "  +0.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.42"
$ws.Range("E39").Value = "  -1.20%  "
$ws.Range("E40").Value = "  -1.40%  "
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("E42").Value = "  +3.64%  "
$ws.Range("E43").Value = "  +0.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.88"
$ws.Range("E44").Value = "  +12.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0486"
$ws.Range("E45").Value = "  +1.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.64"
$ws.Range("E47").Value = "  -2.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.41"
$ws.Range("E48").Value = "  +1.22%  "
$ws.Range("E49").Value = "  +1.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.39"
$ws.Range("E50").Value = "  +5.08%  "
$ws.Range("D51").Value = "2.805.04"
$ws.Range("E51").Value = "  +12.05%  "
